# Weekly update of Fruta / hortaliza data: each existing data row's
# price/date measurements (columns D, J, K, L, M, P) shift down by one
# row (row r gets what row r-1 used to have), row 8 receives the new
# week's observation, and a new last row (124) is appended carrying the
# values that used to live in the last existing row (123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 8
$lastDataRow = 123
$newLastRow = 124

# Columns that actually change values week over week.
$colD = 4   # Fecha
$colJ = 10  # Volumen
$colK = 11  # Precio minimo
$colL = 12  # Precio maximo
$colM = 13  # Precio promedio ponderado
$colP = 16  # Precio $/Kg

# --- Capture the "old" values for rows 8..123 before we overwrite them ---
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldP = @{}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, $colD).Value2
    $oldJ[$r] = $ws.Cells.Item($r, $colJ).Value2
    $oldK[$r] = $ws.Cells.Item($r, $colK).Value2
    $oldL[$r] = $ws.Cells.Item($r, $colL).Value2
    $oldM[$r] = $ws.Cells.Item($r, $colM).Value2
    $oldP[$r] = $ws.Cells.Item($r, $colP).Value2
}

# --- Shift rows 123 down to 9: new row r = old row (r-1) ---
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $prev = $r - 1
    $ws.Cells.Item($r, $colD).Value2 = $oldD[$prev]
    $ws.Cells.Item($r, $colJ).Value2 = $oldJ[$prev]
    $ws.Cells.Item($r, $colK).Value2 = $oldK[$prev]
    $ws.Cells.Item($r, $colL).Value2 = $oldL[$prev]
    $ws.Cells.Item($r, $colM).Value2 = $oldM[$prev]
    $ws.Cells.Item($r, $colP).Value2 = $oldP[$prev]
}

# --- Row 8 gets the brand-new weekly observation ---
$ws.Cells.Item($firstDataRow, $colD).Value2 = 44496
$ws.Cells.Item($firstDataRow, $colJ).Value2 = 2000
# K, L, M, P for row 8 are unchanged (900, 1000, 950, 158)

# --- Append new row 124, a full copy of the old row 123 ---
$ws.Range("A124").Value2 = $ws.Range("A123").Value2
$ws.Range("B124").Value2 = $ws.Range("B123").Value2
$ws.Range("C124").Value2 = $ws.Range("C123").Value2
$ws.Range("D124").Value2 = $oldD[$lastDataRow]
$ws.Range("D124").NumberFormat = $ws.Range("D123").NumberFormat
$ws.Range("E124").Value2 = $ws.Range("E123").Value2
$ws.Range("F124").Value2 = $ws.Range("F123").Value2
$ws.Range("G124").Value2 = $ws.Range("G123").Value2
$ws.Range("H124").Value2 = $ws.Range("H123").Value2
$ws.Range("I124").Value2 = $ws.Range("I123").Value2
$ws.Range("J124").Value2 = $oldJ[$lastDataRow]
$ws.Range("K124").Value2 = $oldK[$lastDataRow]
$ws.Range("L124").Value2 = $oldL[$lastDataRow]
$ws.Range("M124").Value2 = $oldM[$lastDataRow]
$ws.Range("N124").Value2 = $ws.Range("N123").Value2
$ws.Range("O124").Value2 = $ws.Range("O123").Value2
$ws.Range("P124").Value2 = $oldP[$lastDataRow]
$ws.Range("Q124").Value2 = $ws.Range("Q123").Value2
$ws.Range("R124").Value2 = $ws.Range("R123").Value2
